$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Semana 10 de 2024: add rows 11 (week 10, 768 cases) and 12 (week 11, 1 case)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 768

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 1
